# Update essay count from 70+ to 100+ across the client deck (PPTX side).
#
# Touches:
#   - Slide 6: "70以上の思想エッセイ..." summary sentence -> "100以上..."
#   - Slide 8: long bio paragraph "...70以上の思想エッセイを執筆..." -> "...100以上..."
#   - Slide 8: "70+思想エッセイ執筆" badge label -> "100+思想エッセイ執筆"
#   - Slide 8: the badge pill rectangle/label grows a bit wider to fit "100+"
#     instead of "70+", so the two badges to its right shift right by the same
#     amount to stay lined up.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Slide 6 ("3 deliverables") - report description sentence.
# ---------------------------------------------------------------------------
$s6 = $p.Slides.Item(6)
$shp6 = $s6.Shapes.Item(18)
$tr6 = $shp6.TextFrame.TextRange
$tr6.Text = $tr6.Text.Replace("70以上", "100以上")

# ---------------------------------------------------------------------------
# Slide 8 (profile) - long bio paragraph.
# ---------------------------------------------------------------------------
$s8 = $p.Slides.Item(8)
$bioShape = $s8.Shapes.Item(6)
$bioRange = $bioShape.TextFrame.TextRange
$bioRange.Text = $bioRange.Text.Replace("70以上", "100以上")

# ---------------------------------------------------------------------------
# Slide 8 - "70+思想エッセイ執筆" badge (background pill + label) and the
# two badges following it, which shift right to accommodate the wider label.
# ---------------------------------------------------------------------------

# Badge background pill behind the "70+..." label - only grows wider.
$badgeBg = $s8.Shapes.Item(15)
$badgeBg.Width = 97.1999282836914

# "70+思想エッセイ執筆" -> "100+思想エッセイ執筆" label, also grows wider.
$badgeLabel = $s8.Shapes.Item(16)
$badgeLabelRange = $badgeLabel.TextFrame.TextRange
$badgeLabelRange.Text = $badgeLabelRange.Text.Replace("70+", "100+")
$badgeLabel.Width = 82.7999267578125

# Next badge ("オフグリッド実証済み") background pill - shifts right, same size.
$nextBadgeBg = $s8.Shapes.Item(17)
$nextBadgeBg.Left = 570.5999755859375

# Next badge label text itself - shifts right, same size, text unchanged.
$nextBadgeLabel = $s8.Shapes.Item(18)
$nextBadgeLabel.Left = 577.7999267578125
